$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1. Remove the two "ridge / OLS / lasso" paragraphs that open the document.
#    (They are fully replaced by content that already exists further down
#    the document - the "long term bonds..." / "the slope of the yield
#    curve..." paragraphs - so we simply delete these two paragraphs.)
# --------------------------------------------------------------------------
$d.Paragraphs(1).Range.Delete()   # "The ridge regresion is the go to..."
$d.Paragraphs(1).Range.Delete()   # "So ridge is boss. Lasso is the cutthrouat..."

# --------------------------------------------------------------------------
# 2. Collapse the run of now-leading empty paragraphs (there were 3 of them)
#    down to nothing - the "long term bonds..." paragraph becomes the first
#    paragraph of the document.
# --------------------------------------------------------------------------
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()

# --------------------------------------------------------------------------
# 3. Between "the slope of the yield curve (...)" and the red
#    "Following Bolder et al. [2004] ..." paragraph there used to be 3 empty
#    paragraphs; only 1 remains.
# --------------------------------------------------------------------------
$d.Paragraphs(3).Range.Delete()
$d.Paragraphs(3).Range.Delete()

# --------------------------------------------------------------------------
# 4. Fix the typo "di_erence" -> "difference" in the "slope of the yield
#    curve" paragraph.
# --------------------------------------------------------------------------
$d.Content.Find.Execute("di_erence", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "difference", 2) | Out-Null

# --------------------------------------------------------------------------
# 5. Re-create the run break right before "a linear interpolation..." in the
#    last paragraph (cosmetic - matches how the author's edit split the
#    trailing run while leaving the rest of the paragraph untouched).
# --------------------------------------------------------------------------
$tailRng = $d.Content.Duplicate
$tailRng.Find.Execute(" a linear interpolation between the 2-year and 30-year yields.", `
                       $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($tailRng.Find.Found) {
    $tailStart = $tailRng.Start
    $tailText = $tailRng.Text
    $tailRng.Delete()
    $insPt = $d.Range($tailStart, $tailStart)
    $insPt.InsertAfter($tailText)
    $newRng = $d.Range($tailStart, $tailStart + $tailText.Length)
    $newRng.Font.Color = 255
}
